# Applies the LOM3204.xlsx restructuring described by the commit diff.
# Net effect (no row insert/delete needed except dropping the trailing
# two rows): rows 10-21 get new label/value pairs and row heights, and the
# old rows 22-23 are removed entirely so the sheet ends at row 21.
#
# Note: a couple of cells are populated via Range.Copy() from a sibling
# cell in the same column instead of a plain Value assignment. This
# sidesteps two quirks of plain Value-assignment in this workbook:
#   1. A literal "01/01/2016" string typed into a General-formatted cell
#      gets auto-parsed into a date serial number + new date style.
#   2. Column B has two overlapping <col> style definitions in this
#      sheet, so a freshly-created (previously empty) B-cell picks up
#      the wrong default style. Copying an existing, correctly-styled
#      B/C cell guarantees the right style id is reused.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos: now shows the professor identification string ---
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# --- Row 13: "Programa resumido:" label with the activation date value ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B8").Copy($ws.Range("B13"))
$ws.Range("C8").Copy($ws.Range("C13"))
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: "Short syllabus:" label only (B/C cleared entirely) ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: "Programa:" label with the professor name value ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: "Syllabus:" label only (no B/C in this row) ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: "Avaliação:" label only (B/C cleared entirely), default row height ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows.Item(17).RowHeight = $ws.Rows.Item(12).RowHeight
$ws.Rows.Item(17).UseStandardHeight = $true

# --- Row 18: "Método:" label with Katia's name value ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B20").Copy($ws.Range("B18"))
$ws.Range("C20").Copy($ws.Range("C18"))
$ws.Range("B18").Value = "5817692 - Katia Cristiane Gandolpho Candioto"
$ws.Range("C18").Value = "5817692 - Katia Cristiane Gandolpho Candioto"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: "Critério:" label with the teaching method text ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B20").Copy($ws.Range("B19"))
$ws.Range("C20").Copy($ws.Range("C19"))
$ws.Range("B19").Value = "Aulas expositivas, trabalhos e aulas práticas. Aulas com softwares para desenho técnico."
$ws.Range("C19").Value = "Aulas expositivas, trabalhos e aulas práticas. Aulas com softwares para desenho técnico."
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20: "Norma de recuperação:" label with the grading average text ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média aritmética das notas de atividades em aula e extra aula."
$ws.Range("C20").Value = "Média aritmética das notas de atividades em aula e extra aula."
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21: "Bibliografia:" label with the no-recuperação text ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Devido às características práticas da disciplina, não será oferecida recuperação"
$ws.Range("C21").Value = "Devido às características práticas da disciplina, não será oferecida recuperação"
$ws.Rows.Item(21).RowHeight = 120

# --- Remove the trailing two rows (old "Norma de recuperação:" and "Bibliografia:" rows) ---
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(22).Delete()
